$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in Ali Akhbari's (row 14) quiz scores
$ws.Range("F14").Value = 100
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 60
$ws.Range("K14").Value = 0

# Update the view/selection to match the saved window state
$ws.Range("H20").Select()
$excel.ActiveWindow.ScrollRow = 7
